$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.998.73"
$ws.Range("E2").Value = "  +1.65%  "

$ws.Range("D3").Value = "'1.651.30"
$ws.Range("E3").Value = "  +2.01%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'213.90"
$ws.Range("E5").Value = "  +1.20%  "

$ws.Range("D6").Value = "'0.527"
$ws.Range("E6").Value = "  +0.80%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "'23.63"
$ws.Range("E8").Value = "  +3.79%  "

$ws.Range("D9").Value = "'0.267"
$ws.Range("E9").Value = "  +1.54%  "

$ws.Range("E10").Value = "  +0.28%  "

$ws.Range("D11").Value = "'0.0872"
$ws.Range("E11").Value = "  -1.63%  "

$ws.Range("D12").Value = "'1.881.73"

$ws.Range("D13").Value = "'1.651.12"
$ws.Range("E13").Value = "  +1.95%  "

$ws.Range("D14").Value = "'4.09"
$ws.Range("E14").Value = "  +1.65%  "

$ws.Range("D15").Value = "'0.565"
$ws.Range("E15").Value = "  +2.55%  "

$ws.Range("D16").Value = "'65.73"
$ws.Range("E16").Value = "  +1.09%  "

$ws.Range("D17").Value = "'27.968.95"
$ws.Range("E17").Value = "  +1.55%  "

$ws.Range("D18").Value = "'233.28"
$ws.Range("E18").Value = "  +1.31%  "

$ws.Range("D19").Value = "'7.70"
$ws.Range("E19").Value = "  +2.16%  "

$ws.Range("D20").Value = "'0.0₃0724"
$ws.Range("E20").Value = "  +0.55%  "

$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("D22").Value = "'10.71"
$ws.Range("E22").Value = "  +5.46%  "

$ws.Range("D23").Value = "'4.40"
$ws.Range("E23").Value = "  +2.76%  "

$ws.Range("D24").Value = "'2.15"

$ws.Range("D25").Value = "'152.36"
$ws.Range("E25").Value = "  +1.85%  "

$ws.Range("D26").Value = "'6.92"
$ws.Range("E26").Value = "  +1.31%  "

$ws.Range("D27").Value = "'15.77"
$ws.Range("E27").Value = "  +1.24%  "

$ws.Range("D28").Value = "'0.111"
$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("E31").Value = "  +0.63%  "

$ws.Range("D32").Value = "'3.36"
$ws.Range("E32").Value = "  +2.89%  "

$ws.Range("D33").Value = "'1.456.05"
$ws.Range("E33").Value = "  +1.01%  "

$ws.Range("E34").Value = "  +1.45%  "

$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "  +2.26%  "

$ws.Range("E36").Value = "  -0.67%  "

$ws.Range("D37").Value = "'0.893"
$ws.Range("E37").Value = "  +3.53%  "

$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("D39").Value = "'0.563"
$ws.Range("E39").Value = "  +0.23%  "

$ws.Range("D40").Value = "'0.920"
$ws.Range("E40").Value = "  -2.11%  "

$ws.Range("D41").Value = "'69.48"
$ws.Range("E41").Value = "  +2.00%  "

$ws.Range("E42").Value = "  +2.58%  "

$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("D44").Value = "'2.45"
$ws.Range("E44").Value = "  +0.17%  "

$ws.Range("E45").Value = "  +1.05%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'1.80"
$ws.Range("E46").Value = "  +6.32%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'5.36"
$ws.Range("E47").Value = "  -0.98%  "

$ws.Range("D48").Value = "'1.791.02"
$ws.Range("E48").Value = "  +1.66%  "

$ws.Range("D49").Value = "'88.90"
$ws.Range("E49").Value = "  +3.03%  "

$ws.Range("E50").Value = "  +0.26%  "

$ws.Range("E51").Value = "  +1.02%  "
